$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = "3u2jo5s3s"
$ws.Range("B66").Value = "8dpl7g7kb"
